$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# coins: 42 -> 43
$ws.Range("E2").Value = 43

# lost: 2 -> 1
$ws.Range("G2").Value = 1

# ratio: 1 -> 2
$ws.Range("H2").Value = 2
